# Word COM-interop script implementing the "joker-troupe" edit described
# by the commit diff: retitles the review, reshuffles/edits the
# "What we like" bullet list, tweaks one "What we don't like" bullet and
# rewrites the trailing bold/italic summary lines.
#
# We use Range.InsertXML (a FlatOpc "pkg:package" payload) instead of plain
# text replacement so that the original run layout -- in particular the
# leading empty <w:r/> run that precedes the text run in every body / list
# paragraph of this document -- is preserved exactly, rather than being
# collapsed the way Find/Replace or plain Range.Text assignment would do.

$wdParagraph = 4

function New-FlatOpcPackage {
    param([string]$BodyInnerXml)

    return "<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?>" +
        "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
        "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
        "<pkg:xmlData>" +
        "<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body>" +
        $BodyInnerXml +
        "</w:body></w:document>" +
        "</pkg:xmlData></pkg:part></pkg:package>"
}

function Set-ParagraphRangeXml {
    param($Range, [string]$BodyInnerXml)
    $Range.InsertXML((New-FlatOpcPackage $BodyInnerXml))
}

# Returns a Range covering the *whole* paragraph (including its end-of
# paragraph mark) that contains the first match of $SearchText.
function Get-ParagraphRangeForText {
    param($Document, [string]$SearchText)

    $r = $Document.Range(0, 0)
    $null = $r.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $null = $r.Expand($wdParagraph)
    return $r
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Main heading title (always the very first paragraph, so we address
#    it positionally -- this also disambiguates the later Find for the
#    second, identical-at-the-time occurrence near the end of the doc).
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleXml = '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Joker Troupe for Free</w:t></w:r></w:p>'
Set-ParagraphRangeXml $titlePara.Range $titleXml

# ---------------------------------------------------------------------
# 2) "What we like" bullet list: insert two new bullets, keep "Exciting
#    bonus features" where it now belongs, and reword the final bullet
#    while dropping the other two old ones.
# ---------------------------------------------------------------------
$likeStart = $d.Range(0, 0)
$null = $likeStart.Find.Execute("Exciting bonus features", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$likeEnd = $d.Range(0, 0)
$null = $likeEnd.Find.Execute("What we don't like", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Range spanning all four existing "What we like" bullets (from the start
# of the "Exciting bonus features" paragraph through the start of the
# "What we don't like" heading paragraph, exclusive of the heading itself).
$likeRange = $d.Range($likeStart.Start, $likeEnd.Start)

$bulletPPr = '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>'
$likeXml = (
    "<w:p>$bulletPPr<w:r/><w:r><w:t>Unique and exciting gameplay features</w:t></w:r></w:p>" +
    "<w:p>$bulletPPr<w:r/><w:r><w:t>Vibrant colors and brilliant graphics</w:t></w:r></w:p>" +
    "<w:p>$bulletPPr<w:r/><w:r><w:t>Exciting bonus features</w:t></w:r></w:p>" +
    "<w:p>$bulletPPr<w:r/><w:r><w:t>Engaging and entertaining gameplay</w:t></w:r></w:p>"
)
Set-ParagraphRangeXml $likeRange $likeXml

# ---------------------------------------------------------------------
# 3) "What we don't like": reword the "No progressive jackpot feature"
#    bullet (the "Limited number of paylines" bullet is unchanged).
# ---------------------------------------------------------------------
$jackpotRange = Get-ParagraphRangeForText $d "No progressive jackpot feature"
$jackpotXml = "<w:p>$bulletPPr<w:r/><w:r><w:t>Not as famous as other Joker-themed slot games</w:t></w:r></w:p>"
Set-ParagraphRangeXml $jackpotRange $jackpotXml

# ---------------------------------------------------------------------
# 4) Trailing bold "title" restatement line.
# ---------------------------------------------------------------------
$boldRange = Get-ParagraphRangeForText $d "Play Joker Troupe Slot for Free - Review 2021"
$boldXml = '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Joker Troupe for Free</w:t></w:r></w:p>'
Set-ParagraphRangeXml $boldRange $boldXml

# ---------------------------------------------------------------------
# 5) Trailing italic meta-description line.
# ---------------------------------------------------------------------
$italicRange = Get-ParagraphRangeForText $d "Read our unbiased review of Joker Troupe slot game"
$italicXml = '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Joker Troupe and play this exciting slot game for free.</w:t></w:r></w:p>'
Set-ParagraphRangeXml $italicRange $italicXml

Write-Output "done"
